$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")
$ws.Rows(1).Delete()
$ws.Rows(6).Delete()
$ws.Activate()
[void]$ws.Range("A4").Select()
